$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 1494.5
$ws.Range("I13").Value = 1000
$ws.Range("J13").Value = 1989
$ws.Range("K13").Value = 1000
$ws.Range("L13").Value = 1989
$ws.Range("M13").Value = -831
$ws.Range("N13").Value = -2327
$ws.Range("H19").Value = 666.4286
$ws.Range("I19").Value = 606.8
$ws.Range("J19").Value = 815.5
$ws.Range("K19").Value = 606.8
$ws.Range("L19").Value = 815.5
$ws.Range("M19").Value = -431.8
$ws.Range("N19").Value = -1165.5
$ws.Range("H31").Value = 111112800
$ws.Range("I31").Value = 111112800
$ws.Range("K31").Value = 333338400
$ws.Range("M31").Value = -333338170
$ws.Range("H100").Value = 1590.4
$ws.Range("I100").Value = 1414.8572
$ws.Range("J100").Value = 2000
$ws.Range("K100").Value = 1414.8572
$ws.Range("L100").Value = 2000
$ws.Range("M100").Value = -873.8571999999999
$ws.Range("N100").Value = -3082
$ws.Range("H112").Value = 2543.4546
$ws.Range("J112").Value = 2469.1428
$ws.Range("L112").Value = 7407.428400000001
$ws.Range("N112").Value = -9623.428400000001
$ws.Range("H138").Value = 2136.6743
$ws.Range("I138").Value = 1959.2069
$ws.Range("J138").Value = 2504.2856
$ws.Range("K138").Value = 5877.620699999999
$ws.Range("L138").Value = 7512.8568
$ws.Range("M138").Value = -737.6206999999995
$ws.Range("N138").Value = -17792.8568
$ws.Range("H141").Value = 1565.0435
$ws.Range("I141").Value = 1545.2727
$ws.Range("K141").Value = 4635.8181
$ws.Range("M141").Value = 544.1818999999996

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H41").Value = 1993.6666
$ws.Range("I41").Value = 1993.6666
$ws.Range("K41").Value = 1993.6666
$ws.Range("M41").Value = -1579.6666
$ws.Range("H74").Value = 2733.125
$ws.Range("I74").Value = 2766.4285
$ws.Range("K74").Value = 2766.4285
$ws.Range("M74").Value = -1892.4285
$ws.Range("H77").Value = 2733.125
$ws.Range("I77").Value = 2766.4285
$ws.Range("K77").Value = 13832.1425
$ws.Range("M77").Value = -9464.1425
$ws.Range("H122").Value = 2187.5
$ws.Range("I122").Value = 979.5
$ws.Range("K122").Value = 2938.5
$ws.Range("M122").Value = -488.5
$ws.Range("H124").Value = 69524.75
$ws.Range("J124").Value = 69524.75
$ws.Range("L124").Value = 69524.75
$ws.Range("N124").Value = -79344.75

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 38355.4
$ws.Range("J81").Value = 38355.4
$ws.Range("L81").Value = 38355.4
$ws.Range("N81").Value = -40477.4
$ws.Range("H84").Value = 38355.4
$ws.Range("J84").Value = 38355.4
$ws.Range("L84").Value = 115066.2
$ws.Range("N84").Value = -125674.2
$ws.Range("H130").Value = 99998.55499999999
$ws.Range("J130").Value = 99998.55499999999
$ws.Range("L130").Value = 99998.55499999999
$ws.Range("N130").Value = -110038.555

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1216.76
$ws.Range("I58").Value = 1202.2273
$ws.Range("K58").Value = 1202.2273
$ws.Range("M58").Value = -999.2273
$ws.Range("H122").Value = 2115.6667
$ws.Range("I122").Value = 2115.6667
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 6347.000100000001
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -3897.000100000001
$ws.Range("N122").ClearContents()
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()
$ws.Range("H132").Value = 1945.9375
$ws.Range("I132").Value = 1952.5
$ws.Range("K132").Value = 5857.5
$ws.Range("M132").Value = -3327.5
$ws.Range("H134").Value = 1252.25
$ws.Range("I134").Value = 1252.25
$ws.Range("K134").Value = 3756.75
$ws.Range("M134").Value = -1221.75
$ws.Range("H136").Value = 1216.76
$ws.Range("I136").Value = 1202.2273
$ws.Range("K136").Value = 3606.6819
$ws.Range("M136").Value = -1056.6819
$ws.Range("H141").Value = 161998.78
$ws.Range("J141").Value = 161998.78
$ws.Range("L141").Value = 161998.78
$ws.Range("N141").Value = -172358.78

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H116").Value = 2929
$ws.Range("I116").Value = 2929
$ws.Range("K116").Value = 8787
$ws.Range("M116").Value = -5345
$ws.Range("H122").Value = 1004.5
$ws.Range("I122").Value = 1004
$ws.Range("J122").Value = 1005
$ws.Range("K122").Value = 9036
$ws.Range("L122").Value = 9045
$ws.Range("M122").Value = -6586
$ws.Range("N122").Value = -13945
$ws.Range("H131").Value = 529084.7
$ws.Range("I131").Value = 1737.4
$ws.Range("J131").Value = 717423
$ws.Range("K131").Value = 5212.200000000001
$ws.Range("L131").Value = 2152269
$ws.Range("M131").Value = -172.2000000000007
$ws.Range("N131").Value = -2162349
$ws.Range("H141").Value = 7189.6
$ws.Range("I141").Value = 7189.6
$ws.Range("K141").Value = 21568.8
$ws.Range("M141").Value = -16388.8

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H25").Value = 9996.5
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 9996.5
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 9996.5
$ws.Range("M25").ClearContents()
$ws.Range("N25").Value = -11054.5
$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("M41").ClearContents()
$ws.Range("N41").ClearContents()
$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("M44").ClearContents()
$ws.Range("H113").Value = 1233
$ws.Range("I113").Value = 1100
$ws.Range("J113").Value = 1499
$ws.Range("K113").Value = 1100
$ws.Range("L113").Value = 1499
$ws.Range("M113").Value = 1070
$ws.Range("N113").Value = -5839
$ws.Range("H122").Value = 2763.3333
$ws.Range("I122").Value = 2763.3333
$ws.Range("K122").Value = 8289.999899999999
$ws.Range("M122").Value = -5839.999899999999
$ws.Range("H132").Value = 1252.7778
$ws.Range("I132").Value = 1280.1177
$ws.Range("J132").Value = 788
$ws.Range("K132").Value = 3840.3531
$ws.Range("L132").Value = 2364
$ws.Range("M132").Value = -1310.3531
$ws.Range("N132").Value = -7424

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5914.769
$ws.Range("I7").Value = 3026.818
$ws.Range("K7").Value = 3026.818
$ws.Range("M7").Value = -2914.818
$ws.Range("H22").Value = 2174.6667
$ws.Range("I22").Value = 2099.818
$ws.Range("K22").Value = 2099.818
$ws.Range("M22").Value = -1804.818
$ws.Range("H27").Value = 2174.6667
$ws.Range("I27").Value = 2099.818
$ws.Range("K27").Value = 2099.818
$ws.Range("M27").Value = -1992.818
$ws.Range("H40").Value = 2086.182
$ws.Range("I40").Value = 1766.5
$ws.Range("J40").Value = 3524.75
$ws.Range("K40").Value = 1766.5
$ws.Range("L40").Value = 3524.75
$ws.Range("M40").Value = -1630.5
$ws.Range("N40").Value = -3796.75
$ws.Range("H61").Value = 2964
$ws.Range("I61").Value = 3425
$ws.Range("K61").Value = 3425
$ws.Range("M61").Value = -3223
$ws.Range("H82").Value = 3032.5
$ws.Range("I82").Value = 3032.5
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 3032.5
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -2671.5
$ws.Range("N82").ClearContents()
$ws.Range("H85").Value = 3032.5
$ws.Range("I85").Value = 3032.5
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 3032.5
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -1784.5
$ws.Range("N85").ClearContents()
$ws.Range("H113").Value = 2964
$ws.Range("I113").Value = 3425
$ws.Range("K113").Value = 3425
$ws.Range("M113").Value = -1255
$ws.Range("H122").Value = 6231.8125
$ws.Range("J122").Value = 5942
$ws.Range("L122").Value = 17826
$ws.Range("N122").Value = -22726
$ws.Range("H126").Value = 5914.769
$ws.Range("I126").Value = 3026.818
$ws.Range("K126").Value = 9080.454000000002
$ws.Range("M126").Value = -6610.454000000002
$ws.Range("H132").Value = 2065.7273
$ws.Range("I132").Value = 1771.8
$ws.Range("J132").Value = 5005
$ws.Range("K132").Value = 5315.4
$ws.Range("L132").Value = 15015
$ws.Range("M132").Value = -2785.4
$ws.Range("N132").Value = -20075
$ws.Range("H136").Value = 3183.0454
$ws.Range("I136").Value = 2738.1052
$ws.Range("K136").Value = 8214.3156
$ws.Range("M136").Value = -5664.3156

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H124").Value = 27000
$ws.Range("J124").Value = 27000
$ws.Range("L124").Value = 27000
$ws.Range("N124").Value = -36820
$ws.Range("H132").Value = 2763.65
$ws.Range("I132").Value = 2973.25
$ws.Range("K132").Value = 8919.75
$ws.Range("M132").Value = -6389.75
$ws.Range("H136").Value = 2774.4583
$ws.Range("I136").Value = 2880.3044
$ws.Range("K136").Value = 8640.913199999999
$ws.Range("M136").Value = -6090.913199999999
